# Update the cryptocurrency price/volume table (rows 2-51) with the latest
# scraped values, including a few re-ordered coin rows (12/13, 35/36, 47-49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.233.11"
$ws.Range("E2").Value = "  -0.20%  "
# Row 3
$ws.Range("D3").Value = "1.659.64"
$ws.Range("E3").Value = "  -0.53%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.49%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5279"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
# Row 7
$ws.Range("E7").Value = "  -0.42%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2688"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.29%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06386"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07691"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "
# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.627"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.13%  "
# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.670.14"
$ws.Range("E13").Value = "  +0.07%  "
# Row 14
$ws.Range("D14").Value = "1.886.35"
$ws.Range("E14").Value = "  -0.62%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5655"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.97%  "
# Row 16
$ws.Range("D16").Value = "0.0₅8278"
$ws.Range("E16").Value = "  +2.16%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
# Row 18
$ws.Range("D18").Value = "26.195.43"
$ws.Range("E18").Value = "  -0.41%  "
# Row 19
$ws.Range("E19").Value = "  -0.43%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.699"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.62%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.22%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.43%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.009"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.95%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1208"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.81%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.303"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.90%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.64%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.529"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05657"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.11%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.280"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.503"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.390"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.69%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.583"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.93%  "
# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9539"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.94%  "
# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.798"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.80%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.404"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.04%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5783"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.57%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01607"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.989"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.06%  "
# Row 41
$ws.Range("E41").Value = "  -0.39%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8358"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.67%  "
# Row 43
$ws.Range("D43").Value = "1.030.67"
$ws.Range("E43").Value = "  -4.37%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "
# Row 45
$ws.Range("D45").Value = "1.797.25"
$ws.Range("E45").Value = "  -0.61%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.22%  "
# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05361"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.23%  "
# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.86%  "
# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₈104"
$ws.Range("E49").Value = "  +2.70%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.057"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4342"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.57%  "
